$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$ft = $sec.Footers.Item(1)
$shp = $ft.Range.InlineShapes.Item(1)
$shp.Select()
$sel = $word.Selection
Write-Output "Selection InlineShapes count=$($sel.InlineShapes.Count)"
$shp2 = $sel.InlineShapes.Item(1)
$shp2.Width = 500
Write-Output "Width=$($shp2.Width)"
